$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "Correct_answer"

# Rows 2-7 correspond to "Purple" entries -> answer "l"
$ws.Range("D2:D7").Value = "l"

# Rows 8-101 correspond to "Blue" entries -> answer "s"
$ws.Range("D8:D101").Value = "s"

# Update the view to match target state: scrolled position and new selection
$ws.Application.ActiveWindow.ScrollRow = 82
$ws.Range("D8:D101").Select()
